$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.263.33"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.592.55"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.816.04"
$ws.Range("D13").Value = "1.572.44"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "26.244.20"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "1.414.15"
$ws.Range("E33").Value = "  +5.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "1.728.47"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("E51").Value = "  +0.02%  "
